{"js": "// Update the worksheet date header and the twenty-five two-digit by\n// two-digit multiplication prompts to the new day's values.\n// Each \"old\" string is unique within the document, so a simple\n// search-and-replace per pair is unambiguous and order independent.\nconst replacements = [\n  [\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"],\n  [\"68\u00d794=\", \"59\u00d754=\"],\n  [\"36\u00d736=\", \"64\u00d753=\"],\n  [\"14\u00d715=\", \"22\u00d743=\"],\n  [\"93\u00d753=\", \"17\u00d749=\"],\n  [\"62\u00d725=\", \"40\u00d793=\"],\n  [\"37\u00d760=\", \"61\u00d793=\"],\n  [\"58\u00d779=\", \"83\u00d774=\"],\n  [\"95\u00d734=\", \"50\u00d723=\"],\n  [\"12\u00d732=\", \"53\u00d731=\"],\n  [\"87\u00d717=\", \"33\u00d756=\"],\n  [\"81\u00d774=\", \"11\u00d766=\"],\n  [\"95\u00d759=\", \"79\u00d780=\"],\n  [\"21\u00d778=\", \"99\u00d730=\"],\n  [\"33\u00d740=\", \"50\u00d727=\"],\n  [\"83\u00d737=\", \"78\u00d774=\"],\n  [\"33\u00d725=\", \"86\u00d797=\"],\n  [\"94\u00d719=\", \"45\u00d758=\"],\n  [\"56\u00d782=\", \"46\u00d775=\"],\n  [\"95\u00d786=\", \"75\u00d790=\"],\n  [\"73\u00d748=\", \"27\u00d761=\"],\n  [\"71\u00d748=\", \"95\u00d771=\"],\n  [\"74\u00d794=\", \"33\u00d767=\"],\n  [\"26\u00d739=\", \"56\u00d792=\"],\n  [\"54\u00d731=\", \"25\u00d746=\"],\n  [\"31\u00d752=\", \"63\u00d719=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and the twenty-five two-digit by two-digit\n# multiplication prompts to the new day's values.\n# Every \"old\" string is unique within the document, so a Find/Replace\n# (one occurrence at a time) per pair is unambiguous.\n\n$wdReplaceNone = 0\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-03 Saturday\", \"2025-05-04 Sunday\"),\n    @(\"68\u00d794=\", \"59\u00d754=\"),\n    @(\"36\u00d736=\", \"64\u00d753=\"),\n    @(\"14\u00d715=\", \"22\u00d743=\"),\n    @(\"93\u00d753=\", \"17\u00d749=\"),\n    @(\"62\u00d725=\", \"40\u00d793=\"),\n    @(\"37\u00d760=\", \"61\u00d793=\"),\n    @(\"58\u00d779=\", \"83\u00d774=\"),\n    @(\"95\u00d734=\", \"50\u00d723=\"),\n    @(\"12\u00d732=\", \"53\u00d731=\"),\n    @(\"87\u00d717=\", \"33\u00d756=\"),\n    @(\"81\u00d774=\", \"11\u00d766=\"),\n    @(\"95\u00d759=\", \"79\u00d780=\"),\n    @(\"21\u00d778=\", \"99\u00d730=\"),\n    @(\"33\u00d740=\", \"50\u00d727=\"),\n    @(\"83\u00d737=\", \"78\u00d774=\"),\n    @(\"33\u00d725=\", \"86\u00d797=\"),\n    @(\"94\u00d719=\", \"45\u00d758=\"),\n    @(\"56\u00d782=\", \"46\u00d775=\"),\n    @(\"95\u00d786=\", \"75\u00d790=\"),\n    @(\"73\u00d748=\", \"27\u00d761=\"),\n    @(\"71\u00d748=\", \"95\u00d771=\"),\n    @(\"74\u00d794=\", \"33\u00d767=\"),\n    @(\"26\u00d739=\", \"56\u00d792=\"),\n    @(\"54\u00d731=\", \"25\u00d746=\"),\n    @(\"31\u00d752=\", \"63\u00d719=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $newText,\n        $wdReplaceOne\n    )\n\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n"}
